# daily auto push: 2026-02-20 07:12 UTC
#
# Insert a new daily-log row right before the current "2026/12/29" block
# (row 846), shifting every subsequent row down by one. The sheet is a
# simple date/weekday/time/rank log with no formulas, so this is a plain
# "insert a whole row, then fill its four cells" edit.
#
# Column A holds the date as literal text (e.g. "2026/02/20"), not a real
# date value. Assigning a date-looking string straight to .Value makes
# Excel auto-convert it to a date serial (and stamp a date NumberFormat on
# the cell), so we briefly force the cell to Text format first and clear
# the formatting again afterwards once the literal string is safely in
# place - that keeps the cell a plain, unstyled text cell exactly like its
# neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 846..887 down to 847..888, opening up a blank row 846.
$ws.Rows.Item(846).Insert()

# Force column A to Text so the date-like string isn't reinterpreted as a
# real date/serial number.
$ws.Range("A846").NumberFormat = "@"

$ws.Range("A846").Value = "2026/02/20"
$ws.Range("B846").Value = "金"
$ws.Range("C846").Value = 13
$ws.Range("D846").Value = 201

# Drop the temporary Text formatting so the new row's cells end up with no
# explicit style, matching every other data row in the sheet.
$ws.Range("A846:D846").ClearFormats()
